$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# Insert a new column at N (14): the existing "Event " column shifts to O,
# and the new, blank column lands at N.
$ws.Columns.Item(14).Insert()

# The new column N is the renamed "Event" header; O keeps the old "Event "
# header text, which becomes "Correction ".
$ws.Range("N1").Value = "Event"
$ws.Range("O1").Value = "Correction "

# Rows 2-13: the old "Event" column (now at O) was always blank; fill the
# new N column with the same "nan" placeholder used by the other data
# columns in this sheet.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

Write-Output "done"
